# Update countries & provincias Spain
#
# 1) Update the "last updated" timestamp in A1.
# 2) Brasil overtakes Rusia -> swap their ranking rows (4/5/6 area),
#    refreshing Estados Unidos (row4) and Brasil (row5) with new case
#    numbers, while Rusia's row (row6) keeps its previous, unchanged
#    figures now one row lower.
# 3) Reino Unido (row8) gets updated case numbers.
# 4) Canada (row17), Kazajistan (row57) get updated case numbers.
# 5) Irak overtakes Azerbaiyan -> swap their ranking rows (69/70/71
#    area), refreshing Irak (row69) with new numbers, Luxemburgo
#    (row70, unaffected by the swap) also gets refreshed numbers, while
#    Azerbaiyan's row (row71) keeps its previous, unchanged figures now
#    one row lower.
# 6) Cuba (row89), Principado de Andorra (row120) and Gibraltar (row160)
#    get updated case numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 18:05"

# --- Estados Unidos / Brasil / Rusia reorder (rows 4-6) ----------------
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1650677
$ws.Range("C4").Value = 5583
$ws.Range("D4").Value = 403315
$ws.Range("E4").Value = 1149562
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 153
$ws.Range("H4").Value = 97800

$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 339687
$ws.Range("C5").Value = 8797
$ws.Range("D5").Value = 135430
$ws.Range("E5").Value = 182678
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 531
$ws.Range("H5").Value = 21579

$ws.Range("A6").Value = "Rusia"
$ws.Range("B6").Value = 335882
$ws.Range("C6").Value = 9434
$ws.Range("D6").Value = 107936
$ws.Range("E6").Value = 224558
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 139
$ws.Range("H6").Value = 3388

# --- Reino Unido (row 8) -------------------------------------------------
$ws.Range("B8").Value = 257154
$ws.Range("C8").Value = 2959
$ws.Range("G8").Value = 282
$ws.Range("H8").Value = 36675

# --- Canada (row 17) ------------------------------------------------------
$ws.Range("B17").Value = 82892
$ws.Range("C17").Value = 412
$ws.Range("D17").Value = 42973
$ws.Range("E17").Value = 33642
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 6277

# --- Kazajistan (row 57) ---------------------------------------------------
$ws.Range("D57").Value = 4214
$ws.Range("E57").Value = 3670

# --- Irak / Luxemburgo / Azerbaiyan reorder (rows 69-71) -------------------
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 4272
$ws.Range("C69").Value = 308
$ws.Range("D69").Value = 2585
$ws.Range("E69").Value = 1535
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 5
$ws.Range("H69").Value = 152

$ws.Range("A70").Value = "Luxemburgo"
$ws.Range("B70").Value = 3990
$ws.Range("C70").Value = 9
$ws.Range("D70").Value = 3758
$ws.Range("E70").Value = 123
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 109

$ws.Range("A71").Value = "Azerbaiyan"
$ws.Range("B71").Value = 3982
$ws.Range("C71").Value = 127
$ws.Range("D71").Value = 2506
$ws.Range("E71").Value = 1427
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 49

# --- Cuba (row 89) ----------------------------------------------------------
$ws.Range("B89").Value = 1931
$ws.Range("C89").Value = 15
$ws.Range("D89").Value = 1671
$ws.Range("E89").Value = 179

# --- Principado de Andorra (row 120) -----------------------------------------
$ws.Range("D120").Value = 653
$ws.Range("E120").Value = 58

# --- Gibraltar (row 160) ------------------------------------------------------
$ws.Range("B160").Value = 152
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 5
